# Fruta / hortaliza, semanal
# Insert two new daily price rows for "Navel Late" oranges (Vega Modelo de Temuco)
# above the existing data block (old row 449), shifting the rest of the
# table down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows before row 449 -- this pushes old rows 449:511 down to 451:513
$ws.Rows("449:450").Insert()

# ---- New row 449 ----
$ws.Range("A449").Value = 10
$ws.Range("B449").Value = "Vega Modelo de Temuco"
$ws.Range("C449").Value = "La Araucanía"
$ws.Range("D449").Value = 44491
$ws.Range("E449").Value = 9
$ws.Range("F449").Value = "Fruta"
$ws.Range("G449").Value = 100102
$ws.Range("H449").Value = "Cítricos"
$ws.Range("I449").Value = 100102005
$ws.Range("J449").Value = "Naranja"
$ws.Range("K449").Value = "Navel Late"
$ws.Range("L449").Value = "Primera"
$ws.Range("M449").Value = 200
$ws.Range("N449").Value = 8000
$ws.Range("O449").Value = 9000
$ws.Range("P449").Value = 8500
$ws.Range("Q449").Value = "$/bandeja 15 kilos granel"
$ws.Range("R449").Value = "Región de O'Higgins"
$ws.Range("S449").Value = 567
$ws.Range("T449").Value = 15

# ---- New row 450 ----
$ws.Range("A450").Value = 10
$ws.Range("B450").Value = "Vega Modelo de Temuco"
$ws.Range("C450").Value = "La Araucanía"
$ws.Range("D450").Value = 44491
$ws.Range("E450").Value = 9
$ws.Range("F450").Value = "Fruta"
$ws.Range("G450").Value = 100102
$ws.Range("H450").Value = "Cítricos"
$ws.Range("I450").Value = 100102005
$ws.Range("J450").Value = "Naranja"
$ws.Range("K450").Value = "Navel Late"
$ws.Range("L450").Value = "Primera"
$ws.Range("M450").Value = 3
$ws.Range("N450").Value = 180000
$ws.Range("O450").Value = 180000
$ws.Range("P450").Value = 180000
$ws.Range("Q450").Value = "$/bins (400 kilos)"
$ws.Range("R450").Value = "Región de O'Higgins"
$ws.Range("S450").Value = 450
$ws.Range("T450").Value = 400
